# Repull data, push all data, mean calculation
# Update the dSF column (F) values to match re-pulled data

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -7
$ws.Range("F3").Value = 4
$ws.Range("F4").Value = -7
$ws.Range("F5").Value = -3
$ws.Range("F7").Value = 5
$ws.Range("F8").Value = 0
$ws.Range("F9").Value = -4
